$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update definitions that were shortened / reworded ---
$ws.Range("B3").Value = "Institution of origin of the data"
$ws.Range("B11").Value = "CWP grid code or area name when not available"
$ws.Range("B12").Value = "Well-known text of the geometry"

# --- Insert two new rows (LON_CENTROID / LAT_CENTROID) after GEOM_WKT row (row 12) ---
$ws.Rows("13:14").Insert()

$ws.Range("A13").Value = "LON_CENTROID"
$ws.Range("B13").Value = "Longitude (decimal degrees) of the centroid of the fishing ground"

$ws.Range("A14").Value = "LAT_CENTROID"
$ws.Range("B14").Value = "Latitude (decimal degrees) of the centroid of the fishing ground"

# --- Update the capture date / quarter definitions (now rows 15-17) ---
$ws.Range("B15").Value = "Minimum date of capture of the fish (YYYY-MM-DD)"
$ws.Range("B16").Value = "Maximum date of capture of the fish (YYYY-MM-DD)"
$ws.Range("B17").Value = "Quarter of the ""average"" date of capture: 1, 2, 3, 4"

# --- Append two new rows at the bottom (log10FL / log10RD) ---
$ws.Range("A24").Value = "log10FL"
$ws.Range("B24").Value = "Logarithm to base 10 of fork length"

$ws.Range("A25").Value = "log10RD"
$ws.Range("B25").Value = "Logarithm to base 10 of round weight"

# --- Restore the selection that Excel recorded on save ---
$ws.Range("B13").Select()
